# Daily attendance processing - 2025-11-25 07:26:13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 47 (Year 3 / Group A3 / MICROBIOLOGY / session 1) has now been
# recorded (41 of 216 students). Re-use the formatting of an already
# "Recorded" row (green fill) so row 47 switches from the
# "Not Recorded" (pink) style to the "Recorded" (green) style, just
# like every other recorded row in the sheet.
# ------------------------------------------------------------------
$src = $ws.Range("A44:I44")
$dst = $ws.Range("A47:I47")
$src.Copy()
$dst.PasteSpecial(-4122)   # xlPasteFormats: copy formatting only, keep values
$excel.CutCopyMode = $false

$ws.Range('G47').Value = 'nourhan.osama@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('H47').Value = '41/216'
$ws.Range('I47').Value = 'Recorded'

# ------------------------------------------------------------------
# Group Statistics block - Group A3 (row 17): one more session
# recorded, one fewer missing, coverage/attendance % recalculated.
# ------------------------------------------------------------------
$ws.Range('O17').Value = 7
$ws.Range('P17').Value = 0

# Percentage cells are stored as literal text (e.g. "36.8%"), not as
# numbers. A leading apostrophe forces Excel to keep the literal text
# instead of auto-converting it to a numeric percentage; re-applying
# the neighbouring cell's format afterwards keeps the original style
# index instead of Excel forking a new "quote prefix" style variant.
$ws.Range('R17').Value = '''36.8%'
$ws.Range('S17').Value = '''54.6%'
$refQ = $ws.Range("Q17")
$refQ.Copy()
$ws.Range("R17:S17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Overall Class Statistics block (K3:L10): recorded/missing session
# counts and coverage % updated to match the newly recorded session.
# ------------------------------------------------------------------
$ws.Range('L6').Value = 38
$ws.Range('L7').Value = 2
$ws.Range('L9').Value = '''33.3%'
$ws.Range('L10').Value = '''44.3%'
$refL = $ws.Range("L8")
$refL.Copy()
$ws.Range("L9:L10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# "Recorded By" cells: same set of e-mails per session, re-ordered
# (the upstream export no longer guarantees a stable ordering).
# ------------------------------------------------------------------
$ws.Range('G2').Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G18').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G19').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G20').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G21').Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G24').Value = 'lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range('G28').Value = 'marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Range('G37').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G38').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G39').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G40').Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G43').Value = 'menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range('G56').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G57').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G58').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G59').Value = 'nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range('G66').Value = 'Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range('G75').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G76').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G77').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G78').Value = 'nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range('G81').Value = 'lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range('G85').Value = 'Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range('G94').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G95').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G96').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G97').Value = 'nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range('G100').Value = 'menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range('G113').Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range('G115').Value = 'yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
